$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 3848.111
$ws.Range("I80").Value = 1026
$ws.Range("J80").Value = 7375.75
$ws.Range("K80").Value = 3078
$ws.Range("L80").Value = 22127.25
$ws.Range("M80").Value = -2080
$ws.Range("N80").Value = -24123.25

# Row 83
$ws.Range("H83").Value = 3848.111
$ws.Range("I83").Value = 1026
$ws.Range("J83").Value = 7375.75
$ws.Range("K83").Value = 9234
$ws.Range("L83").Value = 66381.75
$ws.Range("M83").Value = -4242
$ws.Range("N83").Value = -76365.75

# Row 137
$ws.Range("H137").Value = 4166.7896
$ws.Range("I137").Value = 7032.6665
$ws.Range("J137").Value = 3629.4375
$ws.Range("K137").Value = 21097.9995
$ws.Range("L137").Value = 10888.3125
$ws.Range("M137").Value = -18547.9995
$ws.Range("N137").Value = -15988.3125

# Row 138
$ws.Range("H138").Value = 3880.5813
$ws.Range("I138").Value = 3890.4
$ws.Range("J138").Value = 3877.606
$ws.Range("K138").Value = 11671.2
$ws.Range("L138").Value = 11632.818
$ws.Range("M138").Value = -6531.200000000001
$ws.Range("N138").Value = -21912.818

$ws = $wb.Worksheets.Item("ARM")
# Row 31
$ws.Range("H31").Value = 4250
$ws.Range("I31").Value = 4250
$ws.Range("K31").Value = 4250
$ws.Range("M31").Value = -3956

# Row 32
$ws.Range("H32").Value = 1502.88
$ws.Range("I32").Value = 1119.7465
$ws.Range("K32").Value = 1119.7465
$ws.Range("M32").Value = -832.7465

# Row 61
$ws.Range("H61").Value = 8309.481
$ws.Range("I61").Value = 7381.263
$ws.Range("J61").Value = 10514
$ws.Range("K61").Value = 7381.263
$ws.Range("L61").Value = 10514
$ws.Range("M61").Value = -7169.263
$ws.Range("N61").Value = -10938

# Row 63
$ws.Range("H63").Value = 1899.1428
$ws.Range("I63").Value = 2250
$ws.Range("J63").Value = 1431.3334
$ws.Range("K63").Value = 2250
$ws.Range("L63").Value = 1431.3334
$ws.Range("M63").Value = -1564
$ws.Range("N63").Value = -2803.3334

# Row 66
$ws.Range("H66").Value = 1899.1428
$ws.Range("I66").Value = 2250
$ws.Range("J66").Value = 1431.3334
$ws.Range("K66").Value = 11250
$ws.Range("L66").Value = 7156.666999999999
$ws.Range("M66").Value = -7818
$ws.Range("N66").Value = -14020.667

# Row 136
$ws.Range("H136").Value = 8309.481
$ws.Range("I136").Value = 7381.263
$ws.Range("J136").Value = 10514
$ws.Range("K136").Value = 22143.789
$ws.Range("L136").Value = 31542
$ws.Range("M136").Value = -19593.789
$ws.Range("N136").Value = -36642

$ws = $wb.Worksheets.Item("BSM")
# Row 102
$ws.Range("H102").Value = 9091.299999999999
$ws.Range("I102").Value = 9213.666999999999
$ws.Range("K102").Value = 9213.666999999999
$ws.Range("M102").Value = -5968.666999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1395.2858
$ws.Range("I22").Value = 386.53845
$ws.Range("J22").Value = 3034.5
$ws.Range("K22").Value = 386.53845
$ws.Range("L22").Value = 3034.5
$ws.Range("M22").Value = -36.53845000000001
$ws.Range("N22").Value = -3734.5

# Row 31
$ws.Range("H31").Value = 38282.355
$ws.Range("I31").Value = 3590.4736
$ws.Range("J31").Value = 93211.164
$ws.Range("K31").Value = 3590.4736
$ws.Range("L31").Value = 93211.164
$ws.Range("M31").Value = -3295.4736
$ws.Range("N31").Value = -93801.164

# Row 34
$ws.Range("H34").Value = 38282.355
$ws.Range("I34").Value = 3590.4736
$ws.Range("J34").Value = 93211.164
$ws.Range("K34").Value = 3590.4736
$ws.Range("L34").Value = 93211.164
$ws.Range("M34").Value = -3388.4736
$ws.Range("N34").Value = -93615.164

# Row 50
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

# Row 51
$ws.Range("H51").Value = 30504.75
$ws.Range("I51").Value = 27710
$ws.Range("J51").Value = 33299.5
$ws.Range("K51").Value = 27710
$ws.Range("L51").Value = 33299.5
$ws.Range("M51").Value = -26974
$ws.Range("N51").Value = -34771.5

# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

# Row 61
$ws.Range("H61").Value = 30504.75
$ws.Range("I61").Value = 27710
$ws.Range("J61").Value = 33299.5
$ws.Range("K61").Value = 27710
$ws.Range("L61").Value = 33299.5
$ws.Range("M61").Value = -27362
$ws.Range("N61").Value = -33995.5

# Row 112
$ws.Range("H112").Value = 75021.336
$ws.Range("J112").Value = 75021.336
$ws.Range("L112").Value = 75021.336
$ws.Range("N112").Value = -77975.336

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 4435.3335
$ws.Range("J68").Value = 4435.3335
$ws.Range("L68").Value = 13306.0005
$ws.Range("N68").Value = -14928.0005

# Row 71
$ws.Range("H71").Value = 4435.3335
$ws.Range("J71").Value = 4435.3335
$ws.Range("L71").Value = 39918.0015
$ws.Range("N71").Value = -48030.0015

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 14937.637
$ws.Range("I70").Value = 5852.0625
$ws.Range("K70").Value = 5852.0625
$ws.Range("M70").Value = -5582.0625

# Row 73
$ws.Range("H73").Value = 14937.637
$ws.Range("I73").Value = 5852.0625
$ws.Range("K73").Value = 5852.0625
$ws.Range("M73").Value = -4916.0625

# Row 80
$ws.Range("H80").Value = 5201.3477
$ws.Range("I80").Value = 4154.778
$ws.Range("J80").Value = 5874.143
$ws.Range("K80").Value = 4154.778
$ws.Range("L80").Value = 5874.143
$ws.Range("M80").Value = -3156.778
$ws.Range("N80").Value = -7870.143

# Row 83
$ws.Range("H83").Value = 5201.3477
$ws.Range("I83").Value = 4154.778
$ws.Range("J83").Value = 5874.143
$ws.Range("K83").Value = 20773.89
$ws.Range("L83").Value = 29370.715
$ws.Range("M83").Value = -15781.89
$ws.Range("N83").Value = -39354.715

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4337.3335
$ws.Range("I22").Value = 1433.4286
$ws.Range("J22").Value = 14501
$ws.Range("K22").Value = 1433.4286
$ws.Range("L22").Value = 14501
$ws.Range("M22").Value = -1138.4286
$ws.Range("N22").Value = -15091

# Row 27
$ws.Range("H27").Value = 4337.3335
$ws.Range("I27").Value = 1433.4286
$ws.Range("J27").Value = 14501
$ws.Range("K27").Value = 1433.4286
$ws.Range("L27").Value = 14501
$ws.Range("M27").Value = -1326.4286
$ws.Range("N27").Value = -14715

# Row 55
$ws.Range("H55").Value = 2587.1667
$ws.Range("I55").Value = 775.5714
$ws.Range("J55").Value = 5123.4
$ws.Range("K55").Value = 775.5714
$ws.Range("L55").Value = 5123.4
$ws.Range("M55").Value = -602.5714
$ws.Range("N55").Value = -5469.4

# Row 122
$ws.Range("H122").Value = 9773.223
$ws.Range("I122").Value = 8619.25
$ws.Range("K122").Value = 25857.75
$ws.Range("M122").Value = -23407.75

$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Range("H61").Value = 30028.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 30028.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 30028.5
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -30612.5

# Row 112
$ws.Range("H112").Value = 40103.332
$ws.Range("J112").Value = 40103.332
$ws.Range("L112").Value = 40103.332
$ws.Range("N112").Value = -43057.332

# Row 126
$ws.Range("H126").Value = 3004
$ws.Range("I126").Value = 1435
$ws.Range("J126").Value = 6665
$ws.Range("K126").Value = 4305
$ws.Range("L126").Value = 19995
$ws.Range("M126").Value = -1835
$ws.Range("N126").Value = -24935

# Row 136
$ws.Range("H136").Value = 2597.7576
$ws.Range("I136").Value = 1557.2
$ws.Range("J136").Value = 13003.333
$ws.Range("K136").Value = 4671.6
$ws.Range("L136").Value = 39009.999
$ws.Range("M136").Value = -2121.6
$ws.Range("N136").Value = -44109.999
